$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '42.987.93'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  -0.17%  '

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.304.75'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  -0.03%  '

$ws.Range("E4").Value = '  -0.06%  '

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '300.02'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -0.19%  '

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '97.82'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '0.511'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -1.64%  '

$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -2.52%  '

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '36.01'
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0790'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.09%  '

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '18.21'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +1.75%  '

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '6.79'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  -1.48%  '

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '2.664.71'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  -0.01%  '

$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '2.313.08'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  -1.11%  '

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '42.927.76'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -0.12%  '

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '12.66'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -5.38%  '

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0903'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -0.55%  '

$ws.Range("E21").Value = '  -1.60%  '

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '67.99'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  -0.56%  '

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '235.89'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  -1.52%  '

$ws.Range("E24").Value = '  -1.43%  '

$ws.Range("E25").Value = '  +0.01%  '

$ws.Range("E26").Value = '  +0.83%  '

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '4.01'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.46%  '

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '25.50'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +3.18%  '

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '165.49'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.22%  '

$ws.Range("E30").Value = '  +0.26%  '

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '9.06'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -1.25%  '

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '33.23'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  -0.26%  '

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '4.85'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  +0.48%  '

$ws.Range("E34").Value = '  +0.01%  '

$ws.Range("E35").Value = '  -4.06%  '

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '16.94'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -6.35%  '

$ws.Range("E37").Value = '  -1.15%  '

$ws.Range("E38").Value = '  -0.46%  '

$ws.Range("E40").Value = '  -1.67%  '

$ws.Range("E41").Value = '  -0.75%  '

$ws.Range("E42").Value = '  -1.77%  '

$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '2.007.77'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  +0.01%  '

$ws.Range("E44").Value = '  -1.79%  '

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '10.08'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '2.10'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -1.84%  '

$ws.Range("E47").Value = '  +1.07%  '

$ws.Range("E48").Value = '  -1.25%  '

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '2.93'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -2.91%  '

$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '53.64'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  -1.39%  '

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '2.532.25'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  +0.11%  '
